# Auto-generated Excel COM-interop edit script
# Updates cryptos list values (prices / volume change %) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.772.67'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '1.719.86'
$ws.Range('E3').Value = '  +0.40%  '
$c = $ws.Range('D4')
$origStyle = $c.Style
$c.Value = "'0.9997"
$c.Style = $origStyle
$ws.Range('E4').Value = '  +0.30%  '
$c = $ws.Range('D5')
$origStyle = $c.Style
$c.Value = "'239.08"
$c.Style = $origStyle
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('E6').Value = '  +0.26%  '
$ws.Range('E7').Value = '  -2.04%  '
$c = $ws.Range('D8')
$origStyle = $c.Style
$c.Value = "'0.2552"
$c.Style = $origStyle
$ws.Range('E8').Value = '  -1.18%  '
$c = $ws.Range('D9')
$origStyle = $c.Style
$c.Value = "'0.06113"
$c.Style = $origStyle
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').Value = '1.718.48'
$ws.Range('E10').Value = '  +0.36%  '
$c = $ws.Range('D11')
$origStyle = $c.Style
$c.Value = "'15.78"
$c.Style = $origStyle
$ws.Range('E11').Value = '  +2.03%  '
$c = $ws.Range('D12')
$origStyle = $c.Style
$c.Value = "'0.06890"
$c.Style = $origStyle
$c = $ws.Range('D13')
$origStyle = $c.Style
$c.Value = "'0.5932"
$c.Style = $origStyle
$ws.Range('E14').Value = '  -1.71%  '
$c = $ws.Range('D15')
$origStyle = $c.Style
$c.Value = "'76.14"
$c.Style = $origStyle
$ws.Range('E15').Value = '  -0.25%  '
$ws.Range('E16').Value = '  +0.33%  '
$ws.Range('D17').Value = '26.685.92'
$ws.Range('E17').Value = '  +1.69%  '
$c = $ws.Range('D18')
$origStyle = $c.Style
$c.Value = "'0.9997"
$c.Style = $origStyle
$ws.Range('E18').Value = '  +0.28%  '
$c = $ws.Range('D19')
$origStyle = $c.Style
$c.Value = "'0.000006998"
$c.Style = $origStyle
$ws.Range('E19').Value = '  -1.33%  '
$c = $ws.Range('D20')
$origStyle = $c.Style
$c.Value = "'11.23"
$c.Style = $origStyle
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').Value = '1.939.46'
$ws.Range('E21').Value = '  +0.26%  '
$c = $ws.Range('D22')
$origStyle = $c.Style
$c.Value = "'4.355"
$c.Style = $origStyle
$ws.Range('E22').Value = '  -0.92%  '
$c = $ws.Range('D23')
$origStyle = $c.Style
$c.Value = "'8.311"
$c.Style = $origStyle
$ws.Range('E23').Value = '  -1.12%  '
$c = $ws.Range('D24')
$origStyle = $c.Style
$c.Value = "'5.050"
$c.Style = $origStyle
$ws.Range('E24').Value = '  +0.38%  '
$c = $ws.Range('D25')
$origStyle = $c.Style
$c.Value = "'140.59"
$c.Style = $origStyle
$ws.Range('E25').Value = '  +3.38%  '
$c = $ws.Range('D26')
$origStyle = $c.Style
$c.Value = "'15.08"
$c.Style = $origStyle
$ws.Range('E26').Value = '  -0.37%  '
$c = $ws.Range('D27')
$origStyle = $c.Style
$c.Value = "'1.785"
$c.Style = $origStyle
$ws.Range('E27').Value = '  +3.51%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Range('D28')
$origStyle = $c.Style
$c.Value = "'1.377"
$c.Style = $origStyle
$ws.Range('E28').Value = '  -1.20%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D29')
$origStyle = $c.Style
$c.Value = "'105.85"
$c.Style = $origStyle
$ws.Range('E29').Value = '  +0.53%  '
$c = $ws.Range('D30')
$origStyle = $c.Style
$c.Value = "'3.933"
$c.Style = $origStyle
$ws.Range('E30').Value = '  +1.80%  '
$c = $ws.Range('D31')
$origStyle = $c.Style
$c.Value = "'0.07874"
$c.Style = $origStyle
$ws.Range('E31').Value = '  -0.69%  '
$c = $ws.Range('D32')
$origStyle = $c.Style
$c.Value = "'3.620"
$c.Style = $origStyle
$ws.Range('E32').Value = '  +0.52%  '
$c = $ws.Range('D33')
$origStyle = $c.Style
$c.Value = "'0.04619"
$c.Style = $origStyle
$ws.Range('E33').Value = '  +4.38%  '
$c = $ws.Range('D34')
$origStyle = $c.Style
$c.Value = "'2.588"
$c.Style = $origStyle
$ws.Range('E34').Value = '  -0.42%  '
$c = $ws.Range('D35')
$origStyle = $c.Style
$c.Value = "'0.9898"
$c.Style = $origStyle
$c = $ws.Range('D36')
$origStyle = $c.Style
$c.Value = "'0.6085"
$c.Style = $origStyle
$ws.Range('E36').Value = '  -1.29%  '
$c = $ws.Range('D37')
$origStyle = $c.Style
$c.Value = "'0.9161"
$c.Style = $origStyle
$ws.Range('E37').Value = '  -1.56%  '
$c = $ws.Range('D38')
$origStyle = $c.Style
$c.Value = "'2.490"
$c.Style = $origStyle
$ws.Range('E38').Value = '  +5.27%  '
$c = $ws.Range('D39')
$origStyle = $c.Style
$c.Value = "'1.962"
$c.Style = $origStyle
$ws.Range('E39').Value = '  -1.15%  '
$ws.Range('E40').Value = '  +0.33%  '
$c = $ws.Range('D41')
$origStyle = $c.Style
$c.Value = "'5.648"
$c.Style = $origStyle
$ws.Range('E41').Value = '  +4.97%  '
$c = $ws.Range('D42')
$origStyle = $c.Style
$c.Value = "'0.01476"
$c.Style = $origStyle
$ws.Range('E42').Value = '  +0.51%  '
$c = $ws.Range('D43')
$origStyle = $c.Style
$c.Value = "'99.96"
$c.Style = $origStyle
$ws.Range('E43').Value = '  +0.09%  '
$c = $ws.Range('D44')
$origStyle = $c.Style
$c.Value = "'0.3771"
$c.Style = $origStyle
$ws.Range('E44').Value = '  -0.70%  '
$c = $ws.Range('D45')
$origStyle = $c.Style
$c.Value = "'6.668"
$c.Style = $origStyle
$ws.Range('E45').Value = '  -2.55%  '
$c = $ws.Range('D46')
$origStyle = $c.Style
$c.Value = "'0.1140"
$c.Style = $origStyle
$ws.Range('E46').Value = '  -0.67%  '
$c = $ws.Range('D47')
$origStyle = $c.Style
$c.Value = "'0.05337"
$c.Style = $origStyle
$ws.Range('E47').Value = '  -0.09%  '
$c = $ws.Range('D48')
$origStyle = $c.Style
$c.Value = "'7.730"
$c.Style = $origStyle
$ws.Range('E48').Value = '  +0.32%  '
$c = $ws.Range('D49')
$origStyle = $c.Style
$c.Value = "'29.59"
$c.Style = $origStyle
$ws.Range('E49').Value = '  -3.02%  '
$c = $ws.Range('D50')
$origStyle = $c.Style
$c.Value = "'1.227"
$c.Style = $origStyle
$ws.Range('E50').Value = '  +1.58%  '
